# Auto-generated script: update market-price derived columns (H-N)
# across the 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect a refreshed Universalis price snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1461
$ws.Range("J32").Value = 981.3333
$ws.Range("L32").Value = 981.3333
$ws.Range("N32").Value = -1633.3333
$ws.Range("H62").Value = 8438.532999999999
$ws.Range("I62").Value = 6433.222
$ws.Range("K62").Value = 6433.222
$ws.Range("M62").Value = -5809.222
$ws.Range("H65").Value = 8438.532999999999
$ws.Range("I65").Value = 6433.222
$ws.Range("K65").Value = 32166.11
$ws.Range("M65").Value = -29046.11
$ws.Range("H80").Value = 76931510
$ws.Range("I80").Value = 142857800
$ws.Range("J80").Value = 17525.5
$ws.Range("K80").Value = 428573400
$ws.Range("L80").Value = 52576.5
$ws.Range("M80").Value = -428572402
$ws.Range("N80").Value = -54572.5
$ws.Range("H83").Value = 76931510
$ws.Range("I83").Value = 142857800
$ws.Range("J83").Value = 17525.5
$ws.Range("K83").Value = 1285720200
$ws.Range("L83").Value = 157729.5
$ws.Range("M83").Value = -1285715208
$ws.Range("N83").Value = -167713.5
$ws.Range("H101").Value = 528.1818
$ws.Range("I101").Value = 566.1667
$ws.Range("J101").Value = 482.6
$ws.Range("K101").Value = 1698.5001
$ws.Range("L101").Value = 1447.8
$ws.Range("M101").Value = -76.50009999999997
$ws.Range("N101").Value = -4691.8
$ws.Range("H132").Value = 4994.757
$ws.Range("I132").Value = 5405.207
$ws.Range("K132").Value = 16215.621
$ws.Range("M132").Value = -13685.621
$ws.Range("H135").Value = 2433
$ws.Range("I135").Value = 1688.6666
$ws.Range("K135").Value = 15197.9994
$ws.Range("M135").Value = -12662.9994
$ws.Range("H138").Value = 2574.4827
$ws.Range("I138").Value = 2332.1538
$ws.Range("J138").Value = 2677.7705
$ws.Range("K138").Value = 6996.4614
$ws.Range("L138").Value = 8033.3115
$ws.Range("M138").Value = -1856.4614
$ws.Range("N138").Value = -18313.3115

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3137
$ws.Range("I45").Value = 1183
$ws.Range("J45").Value = 8999
$ws.Range("K45").Value = 1183
$ws.Range("L45").Value = 8999
$ws.Range("M45").Value = -806
$ws.Range("N45").Value = -9753
$ws.Range("H59").Value = 5000
$ws.Range("J59").Value = 5000
$ws.Range("L59").Value = 5000
$ws.Range("N59").Value = -6608
$ws.Range("H61").Value = 2242.225
$ws.Range("I61").Value = 2021.2858
$ws.Range("K61").Value = 2021.2858
$ws.Range("M61").Value = -1809.2858
$ws.Range("H74").Value = 6259
$ws.Range("J74").Value = 7998.3335
$ws.Range("L74").Value = 7998.3335
$ws.Range("N74").Value = -9746.333500000001
$ws.Range("H77").Value = 6259
$ws.Range("J77").Value = 7998.3335
$ws.Range("L77").Value = 39991.6675
$ws.Range("N77").Value = -48727.6675
$ws.Range("H122").Value = 6857.7407
$ws.Range("I122").Value = 5757.2856
$ws.Range("K122").Value = 17271.8568
$ws.Range("M122").Value = -14821.8568
$ws.Range("H132").Value = 2485.5676
$ws.Range("I132").Value = 1634.44
$ws.Range("J132").Value = 4258.75
$ws.Range("K132").Value = 4903.32
$ws.Range("L132").Value = 12776.25
$ws.Range("M132").Value = -2373.32
$ws.Range("N132").Value = -17836.25
$ws.Range("H136").Value = 2242.225
$ws.Range("I136").Value = 2021.2858
$ws.Range("K136").Value = 6063.857400000001
$ws.Range("M136").Value = -3513.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1433.8334
$ws.Range("I107").Value = 1272.5
$ws.Range("K107").Value = 1272.5
$ws.Range("M107").Value = 647.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2674.3333
$ws.Range("I16").Value = 2654.2
$ws.Range("J16").Value = 2775
$ws.Range("K16").Value = 2654.2
$ws.Range("L16").Value = 2775
$ws.Range("M16").Value = -2367.2
$ws.Range("N16").Value = -3349
$ws.Range("H31").Value = 3353.5
$ws.Range("I31").Value = 2400.2354
$ws.Range("K31").Value = 2400.2354
$ws.Range("M31").Value = -2105.2354
$ws.Range("H34").Value = 3353.5
$ws.Range("I34").Value = 2400.2354
$ws.Range("K34").Value = 2400.2354
$ws.Range("M34").Value = -2198.2354
$ws.Range("H94").Value = 3709.9092
$ws.Range("J94").Value = 3852.1667
$ws.Range("L94").Value = 3852.1667
$ws.Range("N94").Value = -4754.1667
$ws.Range("H99").Value = 1671.3334
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498
$ws.Range("H113").Value = 2674.3333
$ws.Range("I113").Value = 2654.2
$ws.Range("J113").Value = 2775
$ws.Range("K113").Value = 2654.2
$ws.Range("L113").Value = 2775
$ws.Range("M113").Value = -484.1999999999998
$ws.Range("N113").Value = -7115
$ws.Range("H126").Value = 1671.3334
$ws.Range("I126").Value = 1000
$ws.Range("K126").Value = 3000
$ws.Range("M126").Value = -530
$ws.Range("H132").Value = 1305.4375
$ws.Range("I132").Value = 946.0909
$ws.Range("K132").Value = 2838.2727
$ws.Range("M132").Value = -308.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 132091670
$ws.Range("I4").Value = 103605220
$ws.Range("J4").Value = 199747000
$ws.Range("K4").Value = 310815660
$ws.Range("L4").Value = 599241000
$ws.Range("M4").Value = -310815548
$ws.Range("N4").Value = -599241224
$ws.Range("H12").Value = 34.533333
$ws.Range("J12").Value = 40.9
$ws.Range("L12").Value = 122.7
$ws.Range("N12").Value = -468.7
$ws.Range("H131").Value = 2091.2144
$ws.Range("I131").Value = 1047.125
$ws.Range("J131").Value = 3483.3333
$ws.Range("K131").Value = 3141.375
$ws.Range("L131").Value = 10449.9999
$ws.Range("M131").Value = 1898.625
$ws.Range("N131").Value = -20529.9999
$ws.Range("H136").Value = 1525.1428
$ws.Range("I136").Value = 1525.1428
$ws.Range("K136").Value = 4575.428400000001
$ws.Range("M136").Value = 524.5715999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2106.1904
$ws.Range("I102").Value = 2019.5
$ws.Range("J102").Value = 2626.3333
$ws.Range("K102").Value = 2019.5
$ws.Range("L102").Value = 2626.3333
$ws.Range("M102").Value = -397.5
$ws.Range("N102").Value = -5870.3333
$ws.Range("H113").Value = 2054.8572
$ws.Range("I113").Value = 1571.2858
$ws.Range("J113").Value = 2538.4285
$ws.Range("K113").Value = 1571.2858
$ws.Range("L113").Value = 2538.4285
$ws.Range("M113").Value = 598.7141999999999
$ws.Range("N113").Value = -6878.4285
$ws.Range("H132").Value = 7854
$ws.Range("I132").Value = 4992.6665
$ws.Range("K132").Value = 14977.9995
$ws.Range("M132").Value = -12447.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2454
$ws.Range("I61").Value = 2554.158
$ws.Range("J61").Value = 1502.5
$ws.Range("K61").Value = 2554.158
$ws.Range("L61").Value = 1502.5
$ws.Range("M61").Value = -2352.158
$ws.Range("N61").Value = -1906.5
$ws.Range("H68").Value = 1649
$ws.Range("I68").Value = 1464.3334
$ws.Range("J68").Value = 1759.8
$ws.Range("K68").Value = 1464.3334
$ws.Range("L68").Value = 1759.8
$ws.Range("M68").Value = -715.3334
$ws.Range("N68").Value = -3257.8
$ws.Range("H71").Value = 1649
$ws.Range("I71").Value = 1464.3334
$ws.Range("J71").Value = 1759.8
$ws.Range("K71").Value = 7321.666999999999
$ws.Range("L71").Value = 8799
$ws.Range("M71").Value = -3577.666999999999
$ws.Range("N71").Value = -16287
$ws.Range("H113").Value = 2454
$ws.Range("I113").Value = 2554.158
$ws.Range("J113").Value = 1502.5
$ws.Range("K113").Value = 2554.158
$ws.Range("L113").Value = 1502.5
$ws.Range("M113").Value = -384.1579999999999
$ws.Range("N113").Value = -5842.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2929.5
$ws.Range("I126").Value = 1740
$ws.Range("J126").Value = 5416.636
$ws.Range("K126").Value = 5220
$ws.Range("L126").Value = 16249.908
$ws.Range("M126").Value = -2750
$ws.Range("N126").Value = -21189.908
$ws.Range("H132").Value = 3481.138
$ws.Range("I132").Value = 3507.3157
$ws.Range("J132").Value = 3431.4
$ws.Range("K132").Value = 10521.9471
$ws.Range("L132").Value = 10294.2
$ws.Range("M132").Value = -7991.947100000001
$ws.Range("N132").Value = -15354.2
